$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I18").Value = '%'
$ws.Range("J18").Value = 'Uninterpretable'
$ws.Range("I27").Value = 'sv'
$ws.Range("J27").Value = 'Statement-opinion'
$ws.Range("I35").Value = 'aa'
$ws.Range("J35").Value = 'Agree/Accept'
$ws.Range("I42").Value = 'aa'
$ws.Range("J42").Value = 'Agree/Accept'
$ws.Range("I44").Value = 'aa'
$ws.Range("J44").Value = 'Agree/Accept'
$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I95").Value = 'sd'
$ws.Range("J95").Value = 'Statement-non-opinion'
$ws.Range("I102").Value = 'sd'
$ws.Range("J102").Value = 'Statement-non-opinion'
$ws.Range("I107").Value = 'aa'
$ws.Range("J107").Value = 'Agree/Accept'
$ws.Range("I125").Value = 'aa'
$ws.Range("J125").Value = 'Agree/Accept'
$ws.Range("I126").Value = 'aa'
$ws.Range("J126").Value = 'Agree/Accept'
$ws.Range("I128").Value = 'aa'
$ws.Range("J128").Value = 'Agree/Accept'
$ws.Range("I133").Value = 'aa'
$ws.Range("J133").Value = 'Agree/Accept'
$ws.Range("I140").Value = 'aa'
$ws.Range("J140").Value = 'Agree/Accept'
$ws.Range("I152").Value = 'sv'
$ws.Range("J152").Value = 'Statement-opinion'
$ws.Range("I153").Value = 'ba'
$ws.Range("J153").Value = 'Appreciation'
$ws.Range("I154").Value = 'aa'
$ws.Range("J154").Value = 'Agree/Accept'
$ws.Range("I157").Value = 'sd'
$ws.Range("J157").Value = 'Statement-non-opinion'
$ws.Range("I160").Value = 'sd'
$ws.Range("J160").Value = 'Statement-non-opinion'
$ws.Range("I161").Value = 'sd'
$ws.Range("J161").Value = 'Statement-non-opinion'
$ws.Range("I162").Value = 'b'
$ws.Range("J162").Value = 'Acknowledge (Backchannel)'
$ws.Range("I163").Value = 'aa'
$ws.Range("J163").Value = 'Agree/Accept'
$ws.Range("I177").Value = 'aa'
$ws.Range("J177").Value = 'Agree/Accept'
$ws.Range("I178").Value = 'aa'
$ws.Range("J178").Value = 'Agree/Accept'
$ws.Range("I182").Value = 'sd'
$ws.Range("J182").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I195").Value = 'aa'
$ws.Range("J195").Value = 'Agree/Accept'
$ws.Range("I200").Value = 'sd'
$ws.Range("J200").Value = 'Statement-non-opinion'
$ws.Range("I226").Value = 'ba'
$ws.Range("J226").Value = 'Appreciation'
$ws.Range("I228").Value = 'ba'
$ws.Range("J228").Value = 'Appreciation'
$ws.Range("I229").Value = 'ba'
$ws.Range("J229").Value = 'Appreciation'
$ws.Range("I238").Value = 'qy'
$ws.Range("J238").Value = 'Yes-No-Question'
$ws.Range("I243").Value = 'qy'
$ws.Range("J243").Value = 'Yes-No-Question'
$ws.Range("I248").Value = 'sd'
$ws.Range("J248").Value = 'Statement-non-opinion'
$ws.Range("I254").Value = 'ba'
$ws.Range("J254").Value = 'Appreciation'
$ws.Range("I259").Value = 'sd'
$ws.Range("J259").Value = 'Statement-non-opinion'
$ws.Range("I264").Value = 'sd'
$ws.Range("J264").Value = 'Statement-non-opinion'
$ws.Range("I279").Value = 'aa'
$ws.Range("J279").Value = 'Agree/Accept'
$ws.Range("I295").Value = 'sd'
$ws.Range("J295").Value = 'Statement-non-opinion'
$ws.Range("I296").Value = 'ba'
$ws.Range("J296").Value = 'Appreciation'
$ws.Range("I299").Value = '%'
$ws.Range("J299").Value = 'Uninterpretable'
$ws.Range("I310").Value = 'sd'
$ws.Range("J310").Value = 'Statement-non-opinion'
$ws.Range("I320").Value = '%'
$ws.Range("J320").Value = 'Uninterpretable'
$ws.Range("I350").Value = 'aa'
$ws.Range("J350").Value = 'Agree/Accept'
$ws.Range("I352").Value = 'ba'
$ws.Range("J352").Value = 'Appreciation'
$ws.Range("I357").Value = 'b'
$ws.Range("J357").Value = 'Acknowledge (Backchannel)'
$ws.Range("I363").Value = '%'
$ws.Range("J363").Value = 'Uninterpretable'
$ws.Range("I366").Value = 'ba'
$ws.Range("J366").Value = 'Appreciation'
$ws.Range("I402").Value = 'b'
$ws.Range("J402").Value = 'Acknowledge (Backchannel)'
